# spelling + slight changes to the timeline
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week of 42758/42765 (rows 2-3): Unit Tests Panel -> Neural Net Plugin
$ws.Range("B2").Value = "Neural Net Plugin"
$ws.Range("B3").Value = "Neural Net Plugin"

# Week of 42758/42765 (rows 4-5): Neural Net Plugin -> Plugin-interface specification
$ws.Range("B4").Value = "Plugin-interface specification"
$ws.Range("B5").Value = "Plugin-interface specification"

# Spelling fix: "Graph beutification" -> "Graph beautification" (rows 8-10)
$ws.Range("B8").Value = "Graph beautification"
$ws.Range("B9").Value = "Graph beautification"
$ws.Range("B10").Value = "Graph beautification"

# Spelling fix: "Graph beutification specialization" -> "Graph beautification specialization" (row 11)
$ws.Range("B11").Value = "Graph beautification specialization"

# Update the active selection to match the author's last edit location
$ws.Range("B22").Select()
